# Update cryptocurrency price/volume data per Oct 25 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.467.79'
$ws.Cells.Item(2, 5).Value = '  +1.72%  '
$ws.Cells.Item(3, 4).Value = '1.784.36'
$ws.Cells.Item(3, 5).Value = '  -0.18%  '
$ws.Cells.Item(5, 4).Value = "'222.24"
$ws.Cells.Item(5, 5).Value = '  -1.62%  '
$ws.Cells.Item(6, 5).Value = '  -0.82%  '
$ws.Cells.Item(7, 5).Value = '  -0.17%  '
$ws.Cells.Item(8, 4).Value = "'32.45"
$ws.Cells.Item(8, 5).Value = '  +7.77%  '
$ws.Cells.Item(9, 5).Value = '  +0.39%  '
$ws.Cells.Item(10, 4).Value = "'0.0683"
$ws.Cells.Item(10, 5).Value = '  +2.62%  '
$ws.Cells.Item(11, 5).Value = '  +1.07%  '
$ws.Cells.Item(12, 4).Value = '2.044.14'
$ws.Cells.Item(12, 5).Value = '  +0.00%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.797.06'
$ws.Cells.Item(13, 5).Value = '  +0.36%  '
$ws.Cells.Item(14, 2).Value = 'Chainlink'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(14, 4).Value = "'10.98"
$ws.Cells.Item(14, 5).Value = '  +5.46%  '
$ws.Cells.Item(15, 4).Value = '34.501.79'
$ws.Cells.Item(15, 5).Value = '  +1.73%  '
$ws.Cells.Item(16, 4).Value = "'0.629"
$ws.Cells.Item(16, 5).Value = '  +0.63%  '
$ws.Cells.Item(17, 5).Value = '  +2.17%  '
$ws.Cells.Item(18, 4).Value = "'68.59"
$ws.Cells.Item(18, 5).Value = '  -0.61%  '
$ws.Cells.Item(19, 4).Value = "'253.43"
$ws.Cells.Item(19, 5).Value = '  +0.67%  '
$ws.Cells.Item(20, 5).Value = '  +5.44%  '
$ws.Cells.Item(21, 5).Value = '  -0.19%  '
$ws.Cells.Item(22, 5).Value = '  +1.49%  '
$ws.Cells.Item(23, 5).Value = '  -1.28%  '
$ws.Cells.Item(24, 5).Value = '  -0.04%  '
$ws.Cells.Item(25, 4).Value = "'160.52"
$ws.Cells.Item(25, 5).Value = '  +1.36%  '
$ws.Cells.Item(26, 4).Value = "'16.36"
$ws.Cells.Item(26, 5).Value = '  -0.65%  '
$ws.Cells.Item(27, 5).Value = '  +1.38%  '
$ws.Cells.Item(28, 5).Value = '  -0.40%  '
$ws.Cells.Item(30, 5).Value = '  +0.44%  '
$ws.Cells.Item(31, 4).Value = "'3.76"
$ws.Cells.Item(31, 5).Value = '  -2.08%  '
$ws.Cells.Item(32, 5).Value = '  -0.65%  '
$ws.Cells.Item(33, 4).Value = "'3.55"
$ws.Cells.Item(34, 5).Value = '  +1.04%  '
$ws.Cells.Item(35, 4).Value = '1.430.10'
$ws.Cells.Item(35, 5).Value = '  -4.80%  '
$ws.Cells.Item(36, 5).Value = '  +0.97%  '
$ws.Cells.Item(37, 5).Value = '  -1.18%  '
$ws.Cells.Item(38, 4).Value = "'0.0191"
$ws.Cells.Item(38, 5).Value = '  +2.91%  '
$ws.Cells.Item(39, 4).Value = "'85.08"
$ws.Cells.Item(39, 5).Value = '  +2.10%  '
$ws.Cells.Item(40, 4).Value = "'2.80"
$ws.Cells.Item(40, 5).Value = '  +3.02%  '
$ws.Cells.Item(41, 2).Value = 'ARBITRUM'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(41, 4).Value = "'0.916"
$ws.Cells.Item(41, 5).Value = '  +1.85%  '
$ws.Cells.Item(42, 2).Value = 'HuobiToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(42, 4).Value = "'2.33"
$ws.Cells.Item(42, 5).Value = '  -0.91%  '
$ws.Cells.Item(43, 5).Value = '  +1.99%  '
$ws.Cells.Item(44, 4).Value = "'5.97"
$ws.Cells.Item(44, 5).Value = '  +4.39%  '
$ws.Cells.Item(45, 5).Value = '  -1.37%  '
$ws.Cells.Item(46, 5).Value = '  -5.22%  '
$ws.Cells.Item(47, 4).Value = '1.940.95'
$ws.Cells.Item(47, 5).Value = '  +0.05%  '
$ws.Cells.Item(48, 4).Value = "'12.03"
$ws.Cells.Item(48, 5).Value = '  +2.08%  '
$ws.Cells.Item(49, 4).Value = "'104.13"
$ws.Cells.Item(49, 5).Value = '  +6.36%  '
$ws.Cells.Item(50, 4).Value = "'0.999"
$ws.Cells.Item(50, 5).Value = '  -0.31%  '
$ws.Cells.Item(51, 4).Value = "'49.92"
$ws.Cells.Item(51, 5).Value = '  -2.41%  '
